$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row in column A (data starts at row 2)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    if ($ws.Cells.Item($r, 5).Value2 -eq "Retrofitted_2002") {
        $ws.Cells.Item($r, 5).Value = "fullRNASEQ"
    }
    if ($ws.Cells.Item($r, 2).Value2 -eq "Retrofitted_2002") {
        $ws.Cells.Item($r, 2).Value = "H.BROWN"
    }
}

# Restore the selection left by the editor (below the data range)
$ws.Range("D50:G58").Select()
